# Update the "Förändrad" (Changed) date in column C for rows 2-29
# from 45417 (2024-05-05) to 45418 (2024-05-06).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($row = 2; $row -le 29; $row++) {
    $cell = $ws.Cells.Item($row, 3)
    if ($cell.Value2 -eq 45417) {
        $cell.Value2 = 45418
    }
}
